# Revert the connection matrices to the older values (per commit message),
# rename Sheet1 -> exportme, and restore the older view/selection state.

$wb = $excel.ActiveWorkbook

$wsConn = $wb.Worksheets.Item("connections")
$wsExp  = $wb.Worksheets.Item("Sheet1")

# --- Rename sheet ---
$wsExp.Name = "exportme"

# --- Cell value changes on "connections" sheet ---
$wsConn.Range("H15").Value = 1
$wsConn.Range("P17").Value = 1
$wsConn.Range("P18").Value = 0
$wsConn.Range("AC19").Value = 0
$wsConn.Range("V20").Value = 0
$wsConn.Range("H21").Value = 1
$wsConn.Range("Q23").Value = 0
$wsConn.Range("T23").Value = 0
$wsConn.Range("I24").Value = 1
$wsConn.Range("W25").Value = 0
$wsConn.Range("H26").Value = 1
$wsConn.Range("W26").Value = 0
$wsConn.Range("H27").Value = 1
$wsConn.Range("Y27").Value = 0
$wsConn.Range("K30").Value = 1
$wsConn.Range("L31").Value = 1

# --- Cell value changes on "exportme" sheet (formerly Sheet1) ---
$wsExp.Range("A1").Value = ""

$wsExp.Range("O8").Value = 1
$wsExp.Range("T8").Value = 1
$wsExp.Range("AA9").Value = 1
$wsExp.Range("H15").Value = 1
$wsExp.Range("P17").Value = 1
$wsExp.Range("P18").Value = 0
$wsExp.Range("AC19").Value = 0
$wsExp.Range("V20").Value = 0
$wsExp.Range("W20").Value = 0
$wsExp.Range("H23").Value = 0
$wsExp.Range("O23").Value = 0
$wsExp.Range("Q23").Value = 0
$wsExp.Range("T23").Value = 0
$wsExp.Range("U23").Value = 0
$wsExp.Range("H24").Value = 1
$wsExp.Range("W25").Value = 0
$wsExp.Range("H26").Value = 1
$wsExp.Range("I26").Value = 0
$wsExp.Range("W26").Value = 0
$wsExp.Range("H27").Value = 1
$wsExp.Range("Y27").Value = 0
$wsExp.Range("K30").Value = 1
$wsExp.Range("L31").Value = 1

# --- View-state changes ---

# workbook: no sheet should be marked active/selected-tab; connections becomes
# the visually-active sheet again (tabSelected moves there).
$wsConn.Activate()

# connections sheetView: pane moved back to B1, selection moved back to H26
$wsConn.Application.ActiveWindow.ScrollColumn = 2
$wsExp.Range("H26").Select() | Out-Null
$wsConn.Range("H26").Select() | Out-Null

# exportme sheetView: top-left cell Z14, selection W24
$wsExp.Activate()
$wsExp.Application.ActiveWindow.ScrollRow = 14
$wsExp.Application.ActiveWindow.ScrollColumn = 26
$wsExp.Range("W24").Select() | Out-Null

# re-activate connections as the final/selected sheet
$wsConn.Activate()

# column width changes on exportme: single narrower column A, drop col B width
$wsExp.Columns("A").ColumnWidth = 24.140625
$wsExp.Columns("B").ColumnWidth = 8.43
